$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2851.8333
$ws.Range("J17").Value = 2851.8333
$ws.Range("L17").Value = 8555.499899999999
$ws.Range("N17").Value = -8891.499899999999

$ws.Range("H33").Value = 256.07144
$ws.Range("I33").Value = 157
$ws.Range("J33").Value = 388.16666
$ws.Range("K33").Value = 157
$ws.Range("L33").Value = 388.16666
$ws.Range("M33").Value = 72
$ws.Range("N33").Value = -846.16666

$ws.Range("H62").Value = 3261.4211
$ws.Range("I62").Value = 3124.4285
$ws.Range("J62").Value = 3341.3333
$ws.Range("K62").Value = 3124.4285
$ws.Range("L62").Value = 3341.3333
$ws.Range("M62").Value = -2500.4285
$ws.Range("N62").Value = -4589.3333

$ws.Range("H65").Value = 3261.4211
$ws.Range("I65").Value = 3124.4285
$ws.Range("J65").Value = 3341.3333
$ws.Range("K65").Value = 15622.1425
$ws.Range("L65").Value = 16706.6665
$ws.Range("M65").Value = -12502.1425
$ws.Range("N65").Value = -22946.6665

$ws.Range("H69").Value = 2800
$ws.Range("I69").Value = 2800
$ws.Range("K69").Value = 8400
$ws.Range("M69").Value = -7526

$ws.Range("H72").Value = 2800
$ws.Range("I72").Value = 2800
$ws.Range("K72").Value = 25200
$ws.Range("M72").Value = -20832

$ws.Range("H106").Value = 5299
$ws.Range("I106").Value = 5599
$ws.Range("J106").Value = 4999
$ws.Range("K106").Value = 5599
$ws.Range("L106").Value = 4999
$ws.Range("M106").Value = -4968
$ws.Range("N106").Value = -6261

$ws.Range("H137").Value = 2282.5833
$ws.Range("I137").Value = 2198.2
$ws.Range("J137").Value = 2342.8572
$ws.Range("K137").Value = 6594.599999999999
$ws.Range("L137").Value = 7028.571599999999
$ws.Range("M137").Value = -4044.599999999999
$ws.Range("N137").Value = -12128.5716

$ws.Range("H138").Value = 1774.75
$ws.Range("I138").Value = 1249.1666
$ws.Range("K138").Value = 3747.4998
$ws.Range("M138").Value = 1392.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7596.913
$ws.Range("I32").Value = 7234.9443
$ws.Range("J32").Value = 8900
$ws.Range("K32").Value = 7234.9443
$ws.Range("L32").Value = 8900
$ws.Range("M32").Value = -6947.9443
$ws.Range("N32").Value = -9474

$ws.Range("H61").Value = 7333.3335
$ws.Range("I61").Value = 7000
$ws.Range("K61").Value = 7000
$ws.Range("M61").Value = -6788

$ws.Range("H74").Value = 6050
$ws.Range("I74").Value = 5400
$ws.Range("J74").Value = 8000
$ws.Range("K74").Value = 5400
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = -4526
$ws.Range("N74").Value = -9748

$ws.Range("H77").Value = 6050
$ws.Range("I77").Value = 5400
$ws.Range("J77").Value = 8000
$ws.Range("K77").Value = 27000
$ws.Range("L77").Value = 40000
$ws.Range("M77").Value = -22632
$ws.Range("N77").Value = -48736

$ws.Range("H102").Value = 999.3333
$ws.Range("I102").Value = 999.3333
$ws.Range("K102").Value = 999.3333
$ws.Range("M102").Value = 622.6667

$ws.Range("H136").Value = 7333.3335
$ws.Range("I136").Value = 7000
$ws.Range("K136").Value = 21000
$ws.Range("M136").Value = -18450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 998.75
$ws.Range("I107").Value = 785.1429000000001
$ws.Range("K107").Value = 785.1429000000001
$ws.Range("M107").Value = 1134.8571

$ws.Range("H134").Value = 14302.538
$ws.Range("I134").Value = 12994.417
$ws.Range("K134").Value = 38983.251
$ws.Range("M134").Value = -36448.251

$ws.Range("H141").Value = 95000
$ws.Range("J141").Value = 95000
$ws.Range("L141").Value = 95000
$ws.Range("N141").Value = -105360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -650

$ws.Range("H31").Value = 5167.5713
$ws.Range("I31").Value = 4160.6665
$ws.Range("K31").Value = 4160.6665
$ws.Range("M31").Value = -3865.6665

$ws.Range("H34").Value = 5167.5713
$ws.Range("I34").Value = 4160.6665
$ws.Range("K34").Value = 4160.6665
$ws.Range("M34").Value = -3958.6665

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H134").Value = 4229.067
$ws.Range("I134").Value = 4173
$ws.Range("J134").Value = 5014
$ws.Range("K134").Value = 12519
$ws.Range("L134").Value = 15042
$ws.Range("M134").Value = -9984
$ws.Range("N134").Value = -20112

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 682.1667
$ws.Range("I50").Value = 682.1667
$ws.Range("K50").Value = 2046.5001
$ws.Range("M50").Value = -1565.5001

$ws.Range("H53").Value = 682.1667
$ws.Range("I53").Value = 682.1667
$ws.Range("K53").Value = 2046.5001
$ws.Range("M53").Value = -1565.5001

$ws.Range("H131").Value = 1249
$ws.Range("I131").Value = 999
$ws.Range("J131").Value = 1999
$ws.Range("K131").Value = 2997
$ws.Range("L131").Value = 5997
$ws.Range("M131").Value = 2043
$ws.Range("N131").Value = -16077

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 6496.25
$ws.Range("J29").Value = 6496.25
$ws.Range("L29").Value = 6496.25
$ws.Range("N29").Value = -7076.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 17004.5
$ws.Range("I4").Value = 19009
$ws.Range("J4").Value = 15000
$ws.Range("K4").Value = 19009
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = -18896
$ws.Range("N4").Value = -15226

$ws.Range("H16").Value = 747.5
$ws.Range("I16").Value = 747.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 747.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -577.5
$ws.Range("N16").ClearContents()

$ws.Range("H28").Value = 17004.5
$ws.Range("I28").Value = 19009
$ws.Range("J28").Value = 15000
$ws.Range("K28").Value = 19009
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = -18777
$ws.Range("N28").Value = -15464

$ws.Range("H37").Value = 17004.5
$ws.Range("I37").Value = 19009
$ws.Range("J37").Value = 15000
$ws.Range("K37").Value = 19009
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = -18902
$ws.Range("N37").Value = -15214

$ws.Range("H64").Value = 52500
$ws.Range("I64").Value = 75000
$ws.Range("J64").Value = 30000
$ws.Range("K64").Value = 75000
$ws.Range("L64").Value = 30000
$ws.Range("M64").Value = -74775
$ws.Range("N64").Value = -30450

$ws.Range("H67").Value = 52500
$ws.Range("I67").Value = 75000
$ws.Range("J67").Value = 30000
$ws.Range("K67").Value = 75000
$ws.Range("L67").Value = 30000
$ws.Range("M67").Value = -74220
$ws.Range("N67").Value = -31560

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws.Range("H95").Value = 20344
$ws.Range("J95").Value = 20344
$ws.Range("L95").Value = 20344
$ws.Range("N95").Value = -25836

$ws.Range("H132").Value = 11374.5
$ws.Range("I132").Value = 4332.6665
$ws.Range("J132").Value = 32500
$ws.Range("K132").Value = 12997.9995
$ws.Range("L132").Value = 97500
$ws.Range("M132").Value = -10467.9995
$ws.Range("N132").Value = -102560

$ws.Range("H136").Value = 3264.375
$ws.Range("I136").Value = 3264.375
$ws.Range("K136").Value = 9793.125
$ws.Range("M136").Value = -7243.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 1000000
$ws.Range("I29").Value = 1000000
$ws.Range("K29").Value = 1000000
$ws.Range("M29").Value = -999710

$ws.Range("H136").Value = 1069.8889
$ws.Range("I136").Value = 737.7857
$ws.Range("J136").Value = 2232.25
$ws.Range("K136").Value = 2213.3571
$ws.Range("L136").Value = 6696.75
$ws.Range("M136").Value = 336.6428999999998
$ws.Range("N136").Value = -11796.75
